$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R2: Register the information of a client ---
$ws.Range("A22").Value = "R2. Register the information of a client"
$ws.Range("B22").Value = "Main"
$ws.Range("C22").Value = "registerClient():void"

$ws.Range("B23").Value = "Business"
$ws.Range("C23").Value = "registerClient(name : String, lastName : String, id : int, phoneNumber : String, email : String) : String"

$ws.Range("C24").Value = "searchClient(id : int) : model.Client"

$ws.Range("B25").Value = "Client"
$ws.Range("C25").Value = "Client(name : String, lastName : String, id : int, phoneNumber : String, email : String)"

# --- R3: Register the information of a seller ---
$ws.Range("A26").Value = "R3. Register the information of a seller"
$ws.Range("B26").Value = "Main"
$ws.Range("C26").Value = "registerSeller() : void"

$ws.Range("B27").Value = "Business"
$ws.Range("C27").Value = "registerSeller(name : String, lastName : String, id : int) : String"

$ws.Range("C28").Value = "searchSeller(id : int) : model.Seller"

$ws.Range("B29").Value = "Seller"
$ws.Range("C29").Value = "Seller(name : String, lastName : String, id : int)"

# --- Formatting: reuse formatting from analogous existing cells via format copy/paste ---
# A21: filled (no alignment) like the rest of column A's requirement-row fill cells
$ws.Range("C21").Copy()
$ws.Range("A21").PasteSpecial(-4122)

# B22 ("Main"): centered, unfilled, like the other "Main" class-name cells
$ws.Range("B20").Copy()
$ws.Range("B22").PasteSpecial(-4122)

# B23 ("Business"): filled + centered, like the other "Business" class-name cells
$ws.Range("B15").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View / selection state ---
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select() | Out-Null
